# Auto-generated edit script: update Leve profit-table market data cells
# across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2449.2727
$ws.Range("J19").Value = 598.875
$ws.Range("L19").Value = 598.875
$ws.Range("N19").Value = -948.875
$ws.Range("H32").Value = 837.5
$ws.Range("I32").Value = 837.5
$ws.Range("K32").Value = 837.5
$ws.Range("M32").Value = -511.5
$ws.Range("H40").Value = 1404.5294
$ws.Range("I40").Value = 766.6667
$ws.Range("J40").Value = 1752.4546
$ws.Range("K40").Value = 766.6667
$ws.Range("L40").Value = 1752.4546
$ws.Range("M40").Value = -591.6667
$ws.Range("N40").Value = -2102.4546
$ws.Range("H64").Value = 2851.6086
$ws.Range("I64").Value = 2614.6667
$ws.Range("J64").Value = 2935.2354
$ws.Range("K64").Value = 2614.6667
$ws.Range("L64").Value = 2935.2354
$ws.Range("M64").Value = -2366.6667
$ws.Range("N64").Value = -3431.2354
$ws.Range("H67").Value = 2851.6086
$ws.Range("I67").Value = 2614.6667
$ws.Range("J67").Value = 2935.2354
$ws.Range("K67").Value = 2614.6667
$ws.Range("L67").Value = 2935.2354
$ws.Range("M67").Value = -1756.6667
$ws.Range("N67").Value = -4651.2354
$ws.Range("H94").Value = 2800.8
$ws.Range("I94").Value = 2800.8
$ws.Range("K94").Value = 2800.8
$ws.Range("M94").Value = -2349.8
$ws.Range("H129").Value = 1130.2683
$ws.Range("J129").Value = 1261.7428
$ws.Range("L129").Value = 3785.2284
$ws.Range("N129").Value = -13785.2284
$ws.Range("H137").Value = 1925.4736
$ws.Range("I137").Value = 1720.375
$ws.Range("J137").Value = 3019.3333
$ws.Range("K137").Value = 5161.125
$ws.Range("L137").Value = 9057.999899999999
$ws.Range("M137").Value = -2611.125
$ws.Range("N137").Value = -14157.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2787.5
$ws.Range("I2").Value = 2471.4285
$ws.Range("K2").Value = 2471.4285
$ws.Range("M2").Value = -2358.4285
$ws.Range("H4").Value = 610
$ws.Range("I4").Value = 150
$ws.Range("J4").Value = 840
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 840
$ws.Range("M4").Value = -34
$ws.Range("N4").Value = -1072
$ws.Range("H32").Value = 5625.354
$ws.Range("I32").Value = 4475.125
$ws.Range("J32").Value = 11376.5
$ws.Range("K32").Value = 4475.125
$ws.Range("L32").Value = 11376.5
$ws.Range("M32").Value = -4188.125
$ws.Range("N32").Value = -11950.5
$ws.Range("H58").Value = 22000
$ws.Range("J58").Value = 22000
$ws.Range("L58").Value = 22000
$ws.Range("N58").Value = -22860
$ws.Range("H102").Value = 1053.3334
$ws.Range("I102").Value = 1110
$ws.Range("J102").Value = 770
$ws.Range("K102").Value = 1110
$ws.Range("L102").Value = 770
$ws.Range("M102").Value = 512
$ws.Range("N102").Value = -4014
$ws.Range("H116").Value = 2787.5
$ws.Range("I116").Value = 2471.4285
$ws.Range("K116").Value = 2471.4285
$ws.Range("M116").Value = -177.4285
$ws.Range("H132").Value = 17882.281
$ws.Range("I132").Value = 1985.6086
$ws.Range("J132").Value = 58507.11
$ws.Range("K132").Value = 5956.825800000001
$ws.Range("L132").Value = 175521.33
$ws.Range("M132").Value = -3426.825800000001
$ws.Range("N132").Value = -180581.33

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2787.5
$ws.Range("I3").Value = 2471.4285
$ws.Range("K3").Value = 2471.4285
$ws.Range("M3").Value = -2357.4285
$ws.Range("H134").Value = 3049.7368
$ws.Range("I134").Value = 3416.5625
$ws.Range("K134").Value = 10249.6875
$ws.Range("M134").Value = -7714.6875
$ws.Range("H137").Value = 50655
$ws.Range("J137").Value = 50655
$ws.Range("L137").Value = 50655
$ws.Range("N137").Value = -60855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3549.818
$ws.Range("I132").Value = 2581.3333
$ws.Range("K132").Value = 7743.999899999999
$ws.Range("M132").Value = -5213.999899999999
$ws.Range("H134").Value = 1236.1428
$ws.Range("I134").Value = 1126.5
$ws.Range("J134").Value = 1455.4286
$ws.Range("K134").Value = 3379.5
$ws.Range("L134").Value = 4366.2858
$ws.Range("M134").Value = -844.5
$ws.Range("N134").Value = -9436.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 5842.1665
$ws.Range("I63").Value = 5499
$ws.Range("J63").Value = 5910.8
$ws.Range("K63").Value = 16497
$ws.Range("L63").Value = 17732.4
$ws.Range("M63").Value = -15748
$ws.Range("N63").Value = -19230.4
$ws.Range("H66").Value = 5842.1665
$ws.Range("I66").Value = 5499
$ws.Range("J66").Value = 5910.8
$ws.Range("K66").Value = 49491
$ws.Range("L66").Value = 53197.2
$ws.Range("M66").Value = -45747
$ws.Range("N66").Value = -60685.2
$ws.Range("H75").Value = 600.875
$ws.Range("I75").Value = 689.25
$ws.Range("J75").Value = 512.5
$ws.Range("K75").Value = 2067.75
$ws.Range("L75").Value = 1537.5
$ws.Range("M75").Value = -1069.75
$ws.Range("N75").Value = -3533.5
$ws.Range("H78").Value = 600.875
$ws.Range("I78").Value = 689.25
$ws.Range("J78").Value = 512.5
$ws.Range("K78").Value = 6203.25
$ws.Range("L78").Value = 4612.5
$ws.Range("M78").Value = -1211.25
$ws.Range("N78").Value = -14596.5
$ws.Range("H107").Value = 3406.5483
$ws.Range("J107").Value = 200.21428
$ws.Range("L107").Value = 600.64284
$ws.Range("N107").Value = -4440.64284
$ws.Range("H114").Value = 1307.9231
$ws.Range("J114").Value = 1968.375
$ws.Range("L114").Value = 5905.125
$ws.Range("N114").Value = -12413.125
$ws.Range("H117").Value = 1166.6666
$ws.Range("I117").Value = 1166.6666
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 3499.9998
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = -57.99980000000005
$ws.Range("H118").Value = 50002136
$ws.Range("I118").Value = 100000250
$ws.Range("J118").Value = 4025
$ws.Range("K118").Value = 300000750
$ws.Range("L118").Value = 12075
$ws.Range("M118").Value = -299999507
$ws.Range("N118").Value = -14561
$ws.Range("H121").Value = 945.25806
$ws.Range("J121").Value = 1103.9131
$ws.Range("L121").Value = 3311.7393
$ws.Range("N121").Value = -5931.7393
$ws.Range("H129").Value = 11298.637
$ws.Range("I129").Value = 903.3333
$ws.Range("J129").Value = 23773
$ws.Range("K129").Value = 2709.9999
$ws.Range("L129").Value = 71319
$ws.Range("M129").Value = 2290.0001
$ws.Range("N129").Value = -81319
$ws.Range("H131").Value = 701.5816
$ws.Range("J131").Value = 722.3626
$ws.Range("L131").Value = 2167.0878
$ws.Range("N131").Value = -12247.0878
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 19951.678
$ws.Range("I132").Value = 1780.0555
$ws.Range("K132").Value = 5340.166499999999
$ws.Range("M132").Value = -2810.166499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 671498.8
$ws.Range("I132").Value = 928445.0600000001
$ws.Range("K132").Value = 2785335.18
$ws.Range("M132").Value = -2782805.18

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3788592.8
$ws.Range("I107").Value = 731.2857
$ws.Range("J107").Value = 9091599
$ws.Range("K107").Value = 2193.8571
$ws.Range("L107").Value = 27274797
$ws.Range("M107").Value = -273.8571000000002
$ws.Range("N107").Value = -27278637
$ws.Range("H132").Value = 1850.6316
$ws.Range("I132").Value = 1142.4
$ws.Range("J132").Value = 2637.5557
$ws.Range("K132").Value = 3427.2
$ws.Range("L132").Value = 7912.6671
$ws.Range("M132").Value = -897.2000000000003
$ws.Range("N132").Value = -12972.6671
